$d = $word.ActiveDocument

# --- Change 1 ---
# "Füge einem Gegner X Schaden für jedes SYMBOL ABBILDUNG zu."
# merge runs " " + "jedes SYMBOL ABBILDUNG " + "zu." -> " jedes SYMBOL ABBILDUNG zu."
$d.Content.Find.Execute(
    " jedes SYMBOL ABBILDUNG zu.", $true, $false, $false, $false, $false,
    $true, 1, $false, " jedes SYMBOL ABBILDUNG zu.", 2) | Out-Null

# --- Change 2 ---
# "Füge einem Gegner X Schaden und Y Schaden für jedes SYMBOL ABBILDUNG zu."
# merge runs "jedes SYMBOL ABBILDUNG zu" + "." -> "jedes SYMBOL ABBILDUNG zu."
$d.Content.Find.Execute(
    "jedes SYMBOL ABBILDUNG zu.", $true, $false, $false, $false, $false,
    $true, 1, $false, "jedes SYMBOL ABBILDUNG zu.", 2) | Out-Null

# --- Change 3 ---
# merge runs "Füge einem Gegner X" + " Leben hinzu." + " Füge einem Gegner " ->
# "Füge einem Gegner X Leben hinzu. Füge einem Gegner "
$d.Content.Find.Execute(
    "Füge einem Gegner X Leben hinzu. Füge einem Gegner ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Füge einem Gegner X Leben hinzu. Füge einem Gegner ", 2) | Out-Null

# --- Change 4 ---
# add bookmarkStart/bookmarkEnd ("_GoBack") at end of "Füge Y Gegnern jeweils X Schaden zu." paragraph
$yGegnernPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd() -eq "Füge Y Gegnern jeweils X Schaden zu.") {
        $yGegnernPara = $d.Paragraphs.Item($i)
        break
    }
}
$ygEnd = $yGegnernPara.Range.End - 1
$d.Bookmarks.Add("_GoBack", $d.Range($ygEnd, $ygEnd)) | Out-Null

# --- Change 5 ---
# "Drawbacks:" - remove proofErr spellStart/spellEnd and merge runs "Drawbacks" + ":" -> "Drawbacks:"
$drawbacksPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd() -eq "Drawbacks:") {
        $drawbacksPara = $d.Paragraphs.Item($i)
        break
    }
}
$dbStart = $drawbacksPara.Range.Start
$d.Range($dbStart, $dbStart).InsertBefore("Z")
$d.Content.Find.Execute(
    "ZDrawbacks:", $true, $false, $false, $false, $false,
    $true, 1, $false, "Drawbacks:", 2) | Out-Null

# --- Change 6 ---
# merge runs "Füge einem Gegner X Schaden zu. " + "Wirf " -> "Füge einem Gegner X Schaden zu. Wirf "
$d.Content.Find.Execute(
    "Füge einem Gegner X Schaden zu. Wirf ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Füge einem Gegner X Schaden zu. Wirf ", 2) | Out-Null

# --- Change 7 ---
# merge runs " ab, um stattdessen " + "X*2" + " Schaden zuzufügen!" -> " ab, um stattdessen X*2 Schaden zuzufügen!"
$d.Content.Find.Execute(
    " ab, um stattdessen X*2 Schaden zuzufügen!", $true, $false, $false, $false, $false,
    $true, 1, $false, " ab, um stattdessen X*2 Schaden zuzufügen!", 2) | Out-Null

# --- Change 8 ---
# merge runs "Tentakel" + ", " + "Taucheranzug" + ", Meerjungfrau, " -> "Tentakel, Taucheranzug, Meerjungfrau, "
$d.Content.Find.Execute(
    "Tentakel, Taucheranzug, Meerjungfrau, ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Tentakel, Taucheranzug, Meerjungfrau, ", 2) | Out-Null

# --- Change 9 ---
# remove bookmarkStart/bookmarkEnd ("_GoBack") from the last paragraph ("Brawler: weiblich, Farben ")
$d.Bookmarks("_GoBack").Delete()

Write-Output "done"
